# Insert a new row at position 189 (shifting rows 189:280 down to 190:281)
# and populate it with the new record's data, matching the rest of the
# sheet's constant columns (A,B,C,E,F,G,H,Q,R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 189:280 down by one, inserting a blank row at 189.
$ws.Rows.Item(189).Insert()

# Fill in the new row 189 with its values.
$ws.Cells.Item(189, 1).Value = 3
$ws.Cells.Item(189, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(189, 3).Value = "Coquimbo"
$ws.Cells.Item(189, 4).Value = 44523
$ws.Cells.Item(189, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(189, 5).Value = 5
$ws.Cells.Item(189, 6).Value = 100112028
$ws.Cells.Item(189, 7).Value = "Sandia"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 160
$ws.Cells.Item(189, 11).Value = 700
$ws.Cells.Item(189, 12).Value = 700
$ws.Cells.Item(189, 13).Value = 700
$ws.Cells.Item(189, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(189, 15).Value = "Perú"
$ws.Cells.Item(189, 16).Value = 700
$ws.Cells.Item(189, 17).Value = 1
$ws.Cells.Item(189, 18).Value = "Hortaliza"
